# Updates cryptos list values (price + volume%) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.414.92"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.444.52"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.31"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.35"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "2.889.61"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "62.392.08"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "2.446.26"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.69"
$ws.Range("E18").Value = "  -5.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.65"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.86"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.19"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.77"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.54"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "639.32"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "0.0₃0940"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "151.75"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.40"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.34"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.36"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.68"
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0498"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("E51").Value = "  -2.22%  "
